$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-01"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 04-01)"

# Add new data point for April (row 5) in the "2022" column (I)
$ws.Range("I5").Value = 1

# Update the running total for the year (row 14) in column I
$ws.Range("I14").Value = 434
